$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2021563342318059
$ws.Range("C2").Value = 0.5336927223719676
$ws.Range("J2").Value = 0.008086253369272238
$ws.Range("P2").Value = 0.1590296495956873
$ws.Range("S2").Value = 0.09703504043126684
$ws.Range("B3").Value = 0.009900990099009901
$ws.Range("C3").Value = 0.009900990099009901
$ws.Range("J3").Value = 0.0297029702970297
$ws.Range("P3").Value = 0.7277227722772277
$ws.Range("S3").Value = 0.2227722772277228
$ws.Range("J4").Value = 0.01851851851851852
$ws.Range("P4").Value = 0.6481481481481481
$ws.Range("S4").Value = 0.3333333333333333
$ws.Range("B6").Value = 0.08629441624365482
$ws.Range("D6").Value = 0.02538071065989848
$ws.Range("F6").Value = 0.04568527918781726
$ws.Range("J6").Value = 0.2233502538071066
$ws.Range("O6").Value = 0.03045685279187817
$ws.Range("Q6").Value = 0.1319796954314721
$ws.Range("R6").Value = 0.07614213197969544
$ws.Range("S6").Value = 0.3807106598984771
$ws.Range("B7").Value = 0.1555555555555556
$ws.Range("D7").Value = 0.01777777777777778
$ws.Range("F7").Value = 0.03111111111111111
$ws.Range("J7").Value = 0.1155555555555556
$ws.Range("O7").Value = 0.01777777777777778
$ws.Range("Q7").Value = 0.16
$ws.Range("R7").Value = 0.04888888888888889
$ws.Range("S7").Value = 0.4533333333333333
$ws.Range("B8").Value = 0.1317254174397031
$ws.Range("D8").Value = 0.01298701298701299
$ws.Range("E8").Value = 0.001855287569573284
$ws.Range("F8").Value = 0.06679035250463822
$ws.Range("J8").Value = 0.1020408163265306
$ws.Range("O8").Value = 0.02226345083487941
$ws.Range("Q8").Value = 0.1836734693877551
$ws.Range("R8").Value = 0.0575139146567718
$ws.Range("S8").Value = 0.4211502782931354
$ws.Range("B9").Value = 0.08620689655172414
$ws.Range("D9").Value = 0.01149425287356322
$ws.Range("F9").Value = 0.04022988505747126
$ws.Range("J9").Value = 0.1091954022988506
$ws.Range("O9").Value = 0.02298850574712644
$ws.Range("Q9").Value = 0.1724137931034483
$ws.Range("R9").Value = 0.04597701149425287
$ws.Range("S9").Value = 0.5114942528735632
$ws.Range("B10").Value = 0.1237785016286645
$ws.Range("D10").Value = 0.02931596091205212
$ws.Range("F10").Value = 0.06026058631921824
$ws.Range("J10").Value = 0.1091205211726384
$ws.Range("O10").Value = 0.02361563517915309
$ws.Range("Q10").Value = 0.2117263843648209
$ws.Range("R10").Value = 0.05618892508143322
$ws.Range("S10").Value = 0.3859934853420195
$ws.Range("G11").Value = 0.1211180124223603
$ws.Range("J11").Value = 0.06521739130434782
$ws.Range("K11").Value = 0.1925465838509317
$ws.Range("L11").Value = 0.5962732919254659
$ws.Range("S11").Value = 0.02484472049689441
$ws.Range("G12").Value = 0.7563451776649747
$ws.Range("J12").Value = 0.16751269035533
$ws.Range("L12").Value = 0.01015228426395939
$ws.Range("S12").Value = 0.06598984771573604
$ws.Range("F13").Value = 0.01818181818181818
$ws.Range("G13").Value = 0.7636363636363637
$ws.Range("J13").Value = 0.1818181818181818
$ws.Range("S13").Value = 0.03636363636363636
$ws.Range("F15").Value = 0.01363636363636364
$ws.Range("H15").Value = 0.1454545454545454
$ws.Range("I15").Value = 0.06818181818181818
$ws.Range("J15").Value = 0.3090909090909091
$ws.Range("K15").Value = 0.06818181818181818
$ws.Range("M15").Value = 0.00909090909090909
$ws.Range("O15").Value = 0.06818181818181818
$ws.Range("S15").Value = 0.3181818181818182
$ws.Range("F16").Value = 0.01659751037344398
$ws.Range("H16").Value = 0.1701244813278008
$ws.Range("I16").Value = 0.07468879668049792
$ws.Range("J16").Value = 0.3941908713692946
$ws.Range("K16").Value = 0.1120331950207469
$ws.Range("M16").Value = 0.03319502074688797
$ws.Range("N16").Value = 0.004149377593360996
$ws.Range("O16").Value = 0.02904564315352697
$ws.Range("S16").Value = 0.1659751037344398
$ws.Range("F17").Value = 0.0155902004454343
$ws.Range("H17").Value = 0.22271714922049
$ws.Range("I17").Value = 0.08685968819599109
$ws.Range("J17").Value = 0.3674832962138085
$ws.Range("K17").Value = 0.08240534521158129
$ws.Range("M17").Value = 0.0200445434298441
$ws.Range("O17").Value = 0.07572383073496659
$ws.Range("S17").Value = 0.1291759465478842
$ws.Range("F18").Value = 0.03007518796992481
$ws.Range("H18").Value = 0.1879699248120301
$ws.Range("I18").Value = 0.1052631578947368
$ws.Range("J18").Value = 0.4060150375939849
$ws.Range("K18").Value = 0.07518796992481203
$ws.Range("M18").Value = 0.01503759398496241
$ws.Range("O18").Value = 0.05263157894736842
$ws.Range("S18").Value = 0.1278195488721804
$ws.Range("F19").Value = 0.01535240753663643
$ws.Range("H19").Value = 0.2351709699930216
$ws.Range("I19").Value = 0.0642009769713887
$ws.Range("J19").Value = 0.3600837404047453
$ws.Range("K19").Value = 0.115840893230984
$ws.Range("M19").Value = 0.02512212142358688
$ws.Range("O19").Value = 0.05163991625959526
$ws.Range("S19").Value = 0.1325889741800419
